$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 'ECs'
$ws.Range("B2").Value = 'Efna5'
$ws.Range("C2").Value = 'Ephb1'
$ws.Range("D2").Value = 'ECs'
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.2708946666666667
$ws.Range("H2").Value = 0.812684
$ws.Range("I2").Value = 0.1616296696421007
$ws.Range("J2").Value = 0.1616296696421007
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.002008666666667
$ws.Range("N2").Value = 12.006026
$ws.Range("O2").Value = 0.4834231243738785
$ws.Range("P2").Value = 0.4834231243738787
$ws.Range("Q2").Value = 1.084122803753778
$ws.Range("R2").Value = 9.757105233784
$ws.Range("S2").Value = 0.07813551988990212
$ws.Range("T2").Value = 0.07813551988990215

# Row 3
$ws.Range("A3").Value = 'ECs'
$ws.Range("B3").Value = 'Efna5'
$ws.Range("C3").Value = 'Ephb1'
$ws.Range("D3").Value = 'MuSCs'
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.2708946666666667
$ws.Range("H3").Value = 0.812684
$ws.Range("I3").Value = 0.1616296696421007
$ws.Range("J3").Value = 0.1616296696421007
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.265473333333333
$ws.Range("N3").Value = 12.79642
$ws.Range("O3").Value = 0.5152483708764571
$ws.Range("P3").Value = 0.5152483708764573
$ws.Range("Q3").Value = 1.155493976808889
$ws.Range("R3").Value = 10.39944579128
$ws.Range("S3").Value = 0.08327942396839233
$ws.Range("T3").Value = 0.08327942396839234

# Row 4
$ws.Range("A4").Value = 'ECs'
$ws.Range("B4").Value = 'Efna5'
$ws.Range("C4").Value = 'Ephb1'
$ws.Range("D4").Value = 'Resolving-Mac'
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.2708946666666667
$ws.Range("H4").Value = 0.812684
$ws.Range("I4").Value = 0.1616296696421007
$ws.Range("J4").Value = 0.1616296696421007
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.010998
$ws.Range("N4").Value = 0.032994
$ws.Range("O4").Value = 0.001328504749664189
$ws.Range("P4").Value = 0.00132850474966419
$ws.Range("Q4").Value = 0.002979299544
$ws.Range("R4").Value = 0.026813695896
$ws.Range("S4").Value = 0.0002147257838061846
$ws.Range("T4").Value = 0.0002147257838061846

# Row 5
$ws.Range("A5").Value = 'FAPs'
$ws.Range("B5").Value = 'Efna5'
$ws.Range("C5").Value = 'Ephb1'
$ws.Range("D5").Value = 'ECs'
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 1.058121333333333
$ws.Range("H5").Value = 3.174364
$ws.Range("I5").Value = 0.6313295261673385
$ws.Range("J5").Value = 0.6313295261673384
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.002008666666667
$ws.Range("N5").Value = 12.006026
$ws.Range("O5").Value = 0.4834231243738785
$ws.Range("P5").Value = 0.4834231243738787
$ws.Range("Q5").Value = 4.23461074638489
$ws.Range("R5").Value = 38.111496717464
$ws.Range("S5").Value = 0.3051992920492951
$ws.Range("T5").Value = 0.3051992920492951

# Row 6
$ws.Range("A6").Value = 'FAPs'
$ws.Range("B6").Value = 'Efna5'
$ws.Range("C6").Value = 'Ephb1'
$ws.Range("D6").Value = 'MuSCs'
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 1.058121333333333
$ws.Range("H6").Value = 3.174364
$ws.Range("I6").Value = 0.6313295261673385
$ws.Range("J6").Value = 0.6313295261673384
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 4.265473333333333
$ws.Range("N6").Value = 12.79642
$ws.Range("O6").Value = 0.5152483708764571
$ws.Range("P6").Value = 0.5152483708764573
$ws.Range("Q6").Value = 4.513388330764445
$ws.Range("R6").Value = 40.62049497688
$ws.Range("S6").Value = 0.3252915098439268
$ws.Range("T6").Value = 0.3252915098439268

# Row 7
$ws.Range("A7").Value = 'FAPs'
$ws.Range("B7").Value = 'Efna5'
$ws.Range("C7").Value = 'Ephb1'
$ws.Range("D7").Value = 'Resolving-Mac'
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 1.058121333333333
$ws.Range("H7").Value = 3.174364
$ws.Range("I7").Value = 0.6313295261673385
$ws.Range("J7").Value = 0.6313295261673384
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.010998
$ws.Range("N7").Value = 0.032994
$ws.Range("O7").Value = 0.001328504749664189
$ws.Range("P7").Value = 0.00132850474966419
$ws.Range("Q7").Value = 0.011637218424
$ws.Range("R7").Value = 0.104734965816
$ws.Range("S7").Value = 0.0008387242741165513
$ws.Range("T7").Value = 0.0008387242741165514

# Row 8
$ws.Range("A8").Value = 'MuSCs'
$ws.Range("B8").Value = 'Efna5'
$ws.Range("C8").Value = 'Ephb1'
$ws.Range("D8").Value = 'ECs'
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.3470046666666667
$ws.Range("H8").Value = 1.041014
$ws.Range("I8").Value = 0.2070408041905609
$ws.Range("J8").Value = 0.2070408041905609
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 4.002008666666667
$ws.Range("N8").Value = 12.006026
$ws.Range("O8").Value = 0.4834231243738785
$ws.Range("P8").Value = 0.4834231243738787
$ws.Range("Q8").Value = 1.388715683373778
$ws.Range("R8").Value = 12.498441150364
$ws.Range("S8").Value = 0.1000883124346813
$ws.Range("T8").Value = 0.1000883124346814

# Row 9
$ws.Range("A9").Value = 'MuSCs'
$ws.Range("B9").Value = 'Efna5'
$ws.Range("C9").Value = 'Ephb1'
$ws.Range("D9").Value = 'MuSCs'
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.3470046666666667
$ws.Range("H9").Value = 1.041014
$ws.Range("I9").Value = 0.2070408041905609
$ws.Range("J9").Value = 0.2070408041905609
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 4.265473333333333
$ws.Range("N9").Value = 12.79642
$ws.Range("O9").Value = 0.5152483708764571
$ws.Range("P9").Value = 0.5152483708764573
$ws.Range("Q9").Value = 1.480139152208889
$ws.Range("R9").Value = 13.32125236988
$ws.Range("S9").Value = 0.1066774370641381
$ws.Range("T9").Value = 0.1066774370641381

# Row 10
$ws.Range("A10").Value = 'MuSCs'
$ws.Range("B10").Value = 'Efna5'
$ws.Range("C10").Value = 'Ephb1'
$ws.Range("D10").Value = 'Resolving-Mac'
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.3470046666666667
$ws.Range("H10").Value = 1.041014
$ws.Range("I10").Value = 0.2070408041905609
$ws.Range("J10").Value = 0.2070408041905609
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.010998
$ws.Range("N10").Value = 0.032994
$ws.Range("O10").Value = 0.001328504749664189
$ws.Range("P10").Value = 0.00132850474966419
$ws.Range("Q10").Value = 0.003816357324
$ws.Range("R10").Value = 0.03434721591600001
$ws.Range("S10").Value = 0.0002750546917414536
$ws.Range("T10").Value = 0.0002750546917414536
